# Correction In Import Data: add a new "Loan" deduction column (K) to the
# monthly deduction import template, matching the existing header/data styling.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell K1 - same bold/red font as the other headers, but left-aligned
# (no wrap/centered alignment), reusing the font already used by s=3 headers.
$ws.Range("K1").Value = "Loan"
$ws.Range("K1").Font.Bold = $true
$ws.Range("K1").Font.Color = 255

# New data cells K2/K3 - copy the formatting already used by the other data
# cells in the table (style used by column A, etc.) before writing values.
$ws.Range("A2").Copy()
$ws.Range("K2:K3").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("K2").Value = 1000
$ws.Range("K3").Value = 1000

# Page setup touched as part of the save (portrait orientation).
$ws.PageSetup.Orientation = 1

# Leave the selection on the newly added cell, like the saved workbook.
$ws.Range("K3").Select()
